$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header date in C1, copying style (bold/border/center) from B1 first
$ws.Range("B1").Copy($ws.Range("C1"))
$ws.Cells.Item(1, 3).Value = "13-01-2023"

# Update rows 2-41: column A labels, column B values, column C values
$ws.Cells.Item(2, 1).Value = "1810 Renta variable"
$ws.Cells.Item(2, 2).Value = 79414.77
$ws.Cells.Item(2, 3).Value = 79431.98
$ws.Cells.Item(3, 1).Value = "1822 Raices Valores Negociables"
$ws.Cells.Item(3, 2).Value = 221482.79
$ws.Cells.Item(3, 3).Value = 221434.51
$ws.Cells.Item(4, 1).Value = "Adcap IOL Acciones Argentina"
$ws.Cells.Item(4, 2).Value = 33856.63
$ws.Cells.Item(4, 3).Value = 33872.42
$ws.Cells.Item(5, 1).Value = "Allaria Acciones"
$ws.Cells.Item(5, 2).Value = 39693.41
$ws.Cells.Item(5, 3).Value = 39548.23
$ws.Cells.Item(6, 1).Value = "Alpha Acciones"
$ws.Cells.Item(6, 2).Value = 102320.15
$ws.Cells.Item(6, 3).Value = 102067.37
$ws.Cells.Item(7, 1).Value = "Alpha Mega"
$ws.Cells.Item(7, 2).Value = 40361.7
$ws.Cells.Item(7, 3).Value = 40122.64
$ws.Cells.Item(8, 1).Value = "Alpha planeam equil"
$ws.Cells.Item(8, 2).Value = 926.46
$ws.Cells.Item(8, 3).Value = 919.9
$ws.Cells.Item(9, 1).Value = "Argenfunds"
$ws.Cells.Item(9, 2).Value = 12914.81
$ws.Cells.Item(9, 3).Value = 12915.65
$ws.Cells.Item(10, 1).Value = "Arpenta ex Mercosur"
$ws.Cells.Item(10, 2).Value = 4018.68
$ws.Cells.Item(10, 3).Value = 4015.65
$ws.Cells.Item(11, 1).Value = "Balanz"
$ws.Cells.Item(11, 2).Value = 280311.9
$ws.Cells.Item(11, 3).Value = 280821.21
$ws.Cells.Item(12, 1).Value = "Compass Crecimiento"
$ws.Cells.Item(12, 2).Value = 395513.04
$ws.Cells.Item(12, 3).Value = 397251.87
$ws.Cells.Item(13, 1).Value = "Consultatio Acciones Argentina"
$ws.Cells.Item(13, 2).Value = 15851.39
$ws.Cells.Item(13, 3).Value = 15422.33
$ws.Cells.Item(14, 1).Value = "Consultatio Renta Variable"
$ws.Cells.Item(14, 2).Value = 30315.27
$ws.Cells.Item(14, 3).Value = 30361.5
$ws.Cells.Item(15, 1).Value = "Delta Internacional"
$ws.Cells.Item(15, 2).Value = 1583.9
$ws.Cells.Item(15, 3).Value = 1575.92
$ws.Cells.Item(16, 1).Value = "Delta Latinoamerica"
$ws.Cells.Item(16, 2).Value = 5273.87
$ws.Cells.Item(16, 3).Value = 5270.41
$ws.Cells.Item(17, 1).Value = "Delta Select"
$ws.Cells.Item(17, 2).Value = 22206.95
$ws.Cells.Item(17, 3).Value = 22587.44
$ws.Cells.Item(18, 1).Value = "FBA Acciones Argentinas"
$ws.Cells.Item(18, 2).Value = 144020.68
$ws.Cells.Item(18, 3).Value = 149068.41
$ws.Cells.Item(19, 1).Value = "FBA Calificado"
$ws.Cells.Item(19, 2).Value = 139005.78
$ws.Cells.Item(19, 3).Value = 145406.95
$ws.Cells.Item(20, 1).Value = "Fima Acciones"
$ws.Cells.Item(20, 2).Value = 229261.36
$ws.Cells.Item(20, 3).Value = 246788.72
$ws.Cells.Item(21, 1).Value = "Fima PB Acciones"
$ws.Cells.Item(21, 2).Value = 110452.16
$ws.Cells.Item(21, 3).Value = 113029.57
$ws.Cells.Item(22, 1).Value = "Goal Acciones Argentinas"
$ws.Cells.Item(22, 2).Value = 484.7
$ws.Cells.Item(22, 3).Value = 544.77
$ws.Cells.Item(23, 1).Value = "Goal acciones plus"
$ws.Cells.Item(23, 2).Value = 612.96
$ws.Cells.Item(23, 3).Value = 613.24
$ws.Cells.Item(24, 1).Value = "HF Acciones Argentinas"
$ws.Cells.Item(24, 2).Value = 41553.33
$ws.Cells.Item(24, 3).Value = 41353.04
$ws.Cells.Item(25, 1).Value = "HF Acciones Lideres"
$ws.Cells.Item(25, 2).Value = 88052.3
$ws.Cells.Item(25, 3).Value = 88282.33
$ws.Cells.Item(26, 1).Value = "IAM Renta Variable"
$ws.Cells.Item(26, 2).Value = 41973.02
$ws.Cells.Item(26, 3).Value = 43881.66
$ws.Cells.Item(27, 1).Value = "IEB Value"
$ws.Cells.Item(27, 2).Value = 9787.309999999999
$ws.Cells.Item(27, 3).Value = 9791.33
$ws.Cells.Item(28, 1).Value = "Lombardi"
$ws.Cells.Item(28, 2).Value = 17694.57
$ws.Cells.Item(28, 3).Value = 17683.51
$ws.Cells.Item(29, 1).Value = "MAF"
$ws.Cells.Item(29, 2).Value = 5379.84
$ws.Cells.Item(29, 3).Value = 5349.13
$ws.Cells.Item(30, 1).Value = "Megainver"
$ws.Cells.Item(30, 2).Value = 25504.97
$ws.Cells.Item(30, 3).Value = 25435.19
$ws.Cells.Item(31, 1).Value = "Pellegrini Acciones"
$ws.Cells.Item(31, 2).Value = 50379.18
$ws.Cells.Item(31, 3).Value = 65350.43
$ws.Cells.Item(32, 1).Value = "Pionero Acciones"
$ws.Cells.Item(32, 2).Value = 74989.28
$ws.Cells.Item(32, 3).Value = 74852.06
$ws.Cells.Item(33, 1).Value = "Quinquela Acciones"
$ws.Cells.Item(33, 2).Value = 84283.78999999999
$ws.Cells.Item(33, 3).Value = 84278.53999999999
$ws.Cells.Item(34, 1).Value = "Rofex 20 Renta Variable"
$ws.Cells.Item(34, 2).Value = 59443.28
$ws.Cells.Item(34, 3).Value = 59509.98
$ws.Cells.Item(35, 1).Value = "Schroeder RV"
$ws.Cells.Item(35, 2).Value = 253901.02
$ws.Cells.Item(35, 3).Value = 254941.33
$ws.Cells.Item(36, 1).Value = "Supefondo RV"
$ws.Cells.Item(36, 2).Value = 10703.51
$ws.Cells.Item(36, 3).Value = 11401.7
$ws.Cells.Item(37, 1).Value = "Supergestion"
$ws.Cells.Item(37, 2).Value = 345735.97
$ws.Cells.Item(37, 3).Value = 294343.25
$ws.Cells.Item(38, 1).Value = "Toronto Trust Multimercado"
$ws.Cells.Item(38, 2).Value = 29963.54
$ws.Cells.Item(38, 3).Value = 29933.52
$ws.Cells.Item(39, 1).Value = "Toronto trust Argy"
$ws.Cells.Item(39, 2).Value = 72719.12
$ws.Cells.Item(39, 3).Value = 72630.53
$ws.Cells.Item(40, 1).Value = "avg"
$ws.Cells.Item(40, 2).Value = 82156.50999999999
$ws.Cells.Item(40, 3).Value = 82160.22
$ws.Cells.Item(41, 1).Value = "total"
$ws.Cells.Item(41, 2).Value = 3121947.39
$ws.Cells.Item(41, 3).Value = 3122088.22
